# Update cryptocurrency price/volume data per the Fri Mar 10 09:29:40 UTC 2023 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '19.936.51'
$ws.Range('E2').Value = '  -8.02%  '
$ws.Range('D3').Value = '1.407.18'
$ws.Range('E3').Value = '  -8.25%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '1.001'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '271.49'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -5.83%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3694'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -5.58%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3060'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -3.19%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '39.00'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -7.77%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.06505'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -9.19%  '
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.9805'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -6.25%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.002'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.315'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -5.45%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.102'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -7.53%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '16.77'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -9.38%  '
$ws.Range('D16').Value = '1.407.13'
$ws.Range('E16').Value = '  -8.43%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001004'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -8.20%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.05707'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -13.37%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '72.52'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -12.56%  '
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.548'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -9.12%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '14.25'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -7.41%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '10.77'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.54%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.271'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -4.39%  '
$ws.Range('D25').Value = '19.956.38'
$ws.Range('E25').Value = '  -7.93%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.199'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -6.45%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '137.35'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -7.03%  '
$ws.Range('E28').Value = '  -9.20%  '
$ws.Range('D29').Value = '1.565.02'
$ws.Range('E29').Value = '  -8.63%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '108.06'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -7.68%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.841'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -20.64%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.222'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -11.05%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.8127'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -15.19%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.07672'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -5.64%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '8.353'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -3.09%  '
$ws.Range('E36').Value = '  -4.28%  '
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.751'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -6.86%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.1940'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -4.05%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.02021'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -8.13%  '
$ws.Range('E41').Value = '  -7.61%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.064'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -10.03%  '
$ws.Range('B43').Value = 'WEMIXTOKEN'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.290'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -11.10%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.5257'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -8.39%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.506'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -5.82%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '12.08'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -6.85%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5078'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -7.32%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.777'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -4.58%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '109.74'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -4.93%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.035'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -9.92%  '
